$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 8000
$ws.Range("J13").Value = 8000
$ws.Range("L13").Value = 8000
$ws.Range("N13").Value = -8338
$ws.Range("H28").Value = 1583.4667
$ws.Range("I28").Value = 1079.3636
$ws.Range("J28").Value = 2969.75
$ws.Range("K28").Value = 1079.3636
$ws.Range("L28").Value = 2969.75
$ws.Range("M28").Value = -594.3635999999999
$ws.Range("N28").Value = -3939.75
$ws.Range("H33").Value = 396.9
$ws.Range("I33").Value = 396.9
$ws.Range("K33").Value = 396.9
$ws.Range("M33").Value = -167.9
$ws.Range("H43").Value = 7167
$ws.Range("J43").Value = 5500
$ws.Range("L43").Value = 5500
$ws.Range("N43").Value = -5638
$ws.Range("H64").Value = 4999.5
$ws.Range("I64").Value = 4999.6665
$ws.Range("K64").Value = 4999.6665
$ws.Range("M64").Value = -4751.6665
$ws.Range("H67").Value = 4999.5
$ws.Range("I67").Value = 4999.6665
$ws.Range("K67").Value = 4999.6665
$ws.Range("M67").Value = -4141.6665
$ws.Range("H100").Value = 2328.2354
$ws.Range("I100").Value = 2213.2856
$ws.Range("J100").Value = 2864.6667
$ws.Range("K100").Value = 2213.2856
$ws.Range("L100").Value = 2864.6667
$ws.Range("M100").Value = -1672.2856
$ws.Range("N100").Value = -3946.6667
$ws.Range("H125").Value = 1011
$ws.Range("I125").Value = 1032
$ws.Range("J125").Value = 1000.5
$ws.Range("K125").Value = 9288
$ws.Range("L125").Value = 9004.5
$ws.Range("M125").Value = -6828
$ws.Range("N125").Value = -13924.5
$ws.Range("H129").Value = 641.875
$ws.Range("I129").Value = 697.5
$ws.Range("J129").Value = 475
$ws.Range("K129").Value = 2092.5
$ws.Range("L129").Value = 1425
$ws.Range("M129").Value = 2907.5
$ws.Range("N129").Value = -11425

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3615.0908
$ws.Range("I32").Value = 3311.0952
$ws.Range("K32").Value = 3311.0952
$ws.Range("M32").Value = -3024.0952
$ws.Range("H98").Value = 15781.667
$ws.Range("J98").Value = 15781.667
$ws.Range("L98").Value = 15781.667
$ws.Range("N98").Value = -21771.667
$ws.Range("H102").Value = 1397.1904
$ws.Range("I102").Value = 1333.7894
$ws.Range("K102").Value = 1333.7894
$ws.Range("M102").Value = 288.2106000000001
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H110").Value = 10664.8
$ws.Range("I110").Value = 10664.8
$ws.Range("K110").Value = 10664.8
$ws.Range("M110").Value = -8619.799999999999
$ws.Range("H113").Value = 60000
$ws.Range("J113").Value = 60000
$ws.Range("L113").Value = 60000
$ws.Range("N113").Value = -68678
$ws.Range("H131").Value = 81500
$ws.Range("J131").Value = 81500
$ws.Range("L131").Value = 81500
$ws.Range("N131").Value = -91580

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2793.3333
$ws.Range("I20").Value = 1842.4445
$ws.Range("K20").Value = 1842.4445
$ws.Range("M20").Value = -1595.4445
$ws.Range("H86").Value = 3877.923
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 3877.923
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 54984.5
$ws.Range("J20").Value = 54984.5
$ws.Range("L20").Value = 54984.5
$ws.Range("N20").Value = -55456.5
$ws.Range("H30").Value = 54984.5
$ws.Range("J30").Value = 54984.5
$ws.Range("L30").Value = 54984.5
$ws.Range("N30").Value = -55166.5
$ws.Range("H99").Value = 3007.0527
$ws.Range("J99").Value = 3191.1428
$ws.Range("L99").Value = 3191.1428
$ws.Range("N99").Value = -6187.1428
$ws.Range("H105").Value = 1516.2727
$ws.Range("I105").Value = 1434.7
$ws.Range("K105").Value = 1434.7
$ws.Range("M105").Value = 312.3
$ws.Range("H107").Value = 1498.8889
$ws.Range("I107").Value = 450
$ws.Range("K107").Value = 450
$ws.Range("M107").Value = 1470
$ws.Range("H122").Value = 4018.8333
$ws.Range("I122").Value = 4387.5713
$ws.Range("J122").Value = 2728.25
$ws.Range("K122").Value = 13162.7139
$ws.Range("L122").Value = 8184.75
$ws.Range("M122").Value = -10712.7139
$ws.Range("N122").Value = -13084.75
$ws.Range("H126").Value = 3007.0527
$ws.Range("J126").Value = 3191.1428
$ws.Range("L126").Value = 9573.428400000001
$ws.Range("N126").Value = -14513.4284
$ws.Range("H128").Value = 54984.5
$ws.Range("J128").Value = 54984.5
$ws.Range("L128").Value = 54984.5
$ws.Range("N128").Value = -64944.5
$ws.Range("H134").Value = 3224.8333
$ws.Range("I134").Value = 3277.182
$ws.Range("J134").Value = 2649
$ws.Range("K134").Value = 9831.545999999998
$ws.Range("L134").Value = 7947
$ws.Range("M134").Value = -7296.545999999998
$ws.Range("N134").Value = -13017

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7439616.5
$ws.Range("I4").Value = 4619568.5
$ws.Range("K4").Value = 13858705.5
$ws.Range("M4").Value = -13858593.5
$ws.Range("H9").Value = 3216.6667
$ws.Range("I9").Value = 2150
$ws.Range("J9").Value = 3750
$ws.Range("K9").Value = 6450
$ws.Range("L9").Value = 11250
$ws.Range("M9").Value = -6226
$ws.Range("N9").Value = -11698
$ws.Range("H12").Value = 295.58334
$ws.Range("I12").Value = 281
$ws.Range("J12").Value = 324.75
$ws.Range("K12").Value = 843
$ws.Range("L12").Value = 974.25
$ws.Range("M12").Value = -670
$ws.Range("N12").Value = -1320.25
$ws.Range("H35").Value = 737.75
$ws.Range("I35").Value = 737.75
$ws.Range("K35").Value = 2213.25
$ws.Range("M35").Value = -1925.25
$ws.Range("H38").Value = 126.666664
$ws.Range("J38").Value = 90
$ws.Range("L38").Value = 270
$ws.Range("N38").Value = -964
$ws.Range("H93").Value = 14881.667
$ws.Range("I93").Value = 7150
$ws.Range("J93").Value = 18747.5
$ws.Range("K93").Value = 21450
$ws.Range("L93").Value = 56242.5
$ws.Range("M93").Value = -19578
$ws.Range("N93").Value = -59986.5
$ws.Range("H113").Value = 2187.9333
$ws.Range("J113").Value = 2078.4614
$ws.Range("L113").Value = 6235.3842
$ws.Range("N113").Value = -10575.3842
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("M119").ClearContents()
$ws.Range("H131").Value = 57030.145
$ws.Range("J131").Value = 196099.75
$ws.Range("L131").Value = 588299.25
$ws.Range("N131").Value = -598379.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6930.1665
$ws.Range("I70").Value = 6937.4
$ws.Range("J70").Value = 6894
$ws.Range("K70").Value = 6937.4
$ws.Range("L70").Value = 6894
$ws.Range("M70").Value = -6667.4
$ws.Range("N70").Value = -7434
$ws.Range("H73").Value = 6930.1665
$ws.Range("I73").Value = 6937.4
$ws.Range("J73").Value = 6894
$ws.Range("K73").Value = 6937.4
$ws.Range("L73").Value = 6894
$ws.Range("M73").Value = -6001.4
$ws.Range("N73").Value = -8766
$ws.Range("H80").Value = 4624.5
$ws.Range("I80").Value = 4999.3335
$ws.Range("J80").Value = 3500
$ws.Range("K80").Value = 4999.3335
$ws.Range("L80").Value = 3500
$ws.Range("M80").Value = -4001.3335
$ws.Range("N80").Value = -5496
$ws.Range("H83").Value = 4624.5
$ws.Range("I83").Value = 4999.3335
$ws.Range("J83").Value = 3500
$ws.Range("K83").Value = 24996.6675
$ws.Range("L83").Value = 17500
$ws.Range("M83").Value = -20004.6675
$ws.Range("N83").Value = -27484
$ws.Range("H113").Value = 2543
$ws.Range("I113").Value = 2160.2
$ws.Range("K113").Value = 2160.2
$ws.Range("M113").Value = 9.800000000000182
$ws.Range("H126").Value = 1885.2858
$ws.Range("J126").Value = 2298
$ws.Range("L126").Value = 6894
$ws.Range("N126").Value = -11834
$ws.Range("H128").Value = 61191.2
$ws.Range("I128").Value = 45999
$ws.Range("J128").Value = 64989.25
$ws.Range("K128").Value = 45999
$ws.Range("L128").Value = 64989.25
$ws.Range("M128").Value = -41019
$ws.Range("N128").Value = -74949.25
$ws.Range("H141").Value = 69000
$ws.Range("J141").Value = 69000
$ws.Range("L141").Value = 69000
$ws.Range("N141").Value = -79360

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6003
$ws.Range("I7").Value = 6153.2144
$ws.Range("K7").Value = 6153.2144
$ws.Range("M7").Value = -6041.2144
$ws.Range("H82").Value = 1615.6
$ws.Range("I82").Value = 1930.6
$ws.Range("J82").Value = 1300.6
$ws.Range("K82").Value = 1930.6
$ws.Range("L82").Value = 1300.6
$ws.Range("M82").Value = -1569.6
$ws.Range("N82").Value = -2022.6
$ws.Range("H85").Value = 1615.6
$ws.Range("I85").Value = 1930.6
$ws.Range("J85").Value = 1300.6
$ws.Range("K85").Value = 1930.6
$ws.Range("L85").Value = 1300.6
$ws.Range("M85").Value = -682.5999999999999
$ws.Range("N85").Value = -3796.6
$ws.Range("H122").Value = 3692.6667
$ws.Range("I122").Value = 3340.6428
$ws.Range("K122").Value = 10021.9284
$ws.Range("M122").Value = -7571.928400000001
$ws.Range("H126").Value = 6003
$ws.Range("I126").Value = 6153.2144
$ws.Range("K126").Value = 18459.6432
$ws.Range("M126").Value = -15989.6432
$ws.Range("H128").Value = 79999.75
$ws.Range("J128").Value = 79999.75
$ws.Range("L128").Value = 79999.75
$ws.Range("N128").Value = -89959.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 15750
$ws.Range("J29").Value = 11468.75
$ws.Range("L29").Value = 11468.75
$ws.Range("N29").Value = -12048.75
$ws.Range("H74").Value = 17903.4
$ws.Range("J74").Value = 16880.5
$ws.Range("L74").Value = 16880.5
$ws.Range("N74").Value = -18752.5
$ws.Range("H77").Value = 17903.4
$ws.Range("J77").Value = 16880.5
$ws.Range("L77").Value = 50641.5
$ws.Range("N77").Value = -60001.5
$ws.Range("H81").Value = 7000
$ws.Range("J81").Value = 7000
$ws.Range("L81").Value = 14000
$ws.Range("N81").Value = -16122
$ws.Range("H84").Value = 7000
$ws.Range("J84").Value = 7000
$ws.Range("L84").Value = 70000
$ws.Range("N84").Value = -80608
$ws.Range("H96").Value = 1799
$ws.Range("I96").Value = 1799
$ws.Range("K96").Value = 1799
$ws.Range("M96").Value = -426
$ws.Range("H117").Value = 67500
$ws.Range("J117").Value = 67500
$ws.Range("L117").Value = 67500
$ws.Range("N117").Value = -76678
$ws.Range("H130").Value = 54666
$ws.Range("J130").Value = 54666
$ws.Range("L130").Value = 54666
$ws.Range("N130").Value = -64706
